# Update loading_percent values for rows 2-25 (columns C-I, L, M)
# per case "380 kV" rerun results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Row 2
$ws.Range("C2").Value = 8.338680776672632
$ws.Range("D2").Value = 5.218470269711807
$ws.Range("E2").Value = 13.24567078061751
$ws.Range("F2").Value = 27.62804914065313
$ws.Range("G2").Value = 34.84086742476234
$ws.Range("H2").Value = 15.54262918657586
$ws.Range("I2").Value = 23.90876985758839
$ws.Range("L2").Value = 9.798203811002125
$ws.Range("M2").Value = 59.49720345784911

# Row 3
$ws.Range("C3").Value = 8.291255343840522
$ws.Range("D3").Value = 5.321346513211649
$ws.Range("E3").Value = 12.91890434566484
$ws.Range("F3").Value = 27.95247333358189
$ws.Range("G3").Value = 35.33738637258423
$ws.Range("H3").Value = 15.75799061344303
$ws.Range("I3").Value = 24.37116412154704
$ws.Range("L3").Value = 9.642830058159744
$ws.Range("M3").Value = 56.11990664958039

# Row 4
$ws.Range("C4").Value = 8.264461894601629
$ws.Range("D4").Value = 5.386888881226631
$ws.Range("E4").Value = 12.719136007257
$ws.Range("F4").Value = 28.17623514624267
$ws.Range("G4").Value = 35.68616849880728
$ws.Range("H4").Value = 15.89878886335917
$ws.Range("I4").Value = 24.67223097884669
$ws.Range("L4").Value = 9.549971847505226
$ws.Range("M4").Value = 53.93329370532977

# Row 5
$ws.Range("C5").Value = 8.254133057862607
$ws.Range("D5").Value = 5.414197553285003
$ws.Range("E5").Value = 12.63806245967591
$ws.Range("F5").Value = 28.273408250555
$ws.Range("G5").Value = 35.83890097609485
$ws.Range("H5").Value = 15.95828536751692
$ws.Range("I5").Value = 24.79914838591204
$ws.Range("L5").Value = 9.512809553990332
$ws.Range("M5").Value = 53.0141374030155

# Row 6
$ws.Range("C6").Value = 8.252453710482213
$ws.Range("D6").Value = 5.418768421224114
$ws.Range("E6").Value = 12.62462386967004
$ws.Range("F6").Value = 28.28990004140089
$ws.Range("G6").Value = 35.86488966937095
$ws.Range("H6").Value = 15.96829181511723
$ws.Range("I6").Value = 24.82047610748661
$ws.Range("L6").Value = 9.506680870410172
$ws.Range("M6").Value = 52.8598261499986

# Row 7
$ws.Range("C7").Value = 8.26432020256448
$ws.Range("D7").Value = 5.387254744605055
$ws.Range("E7").Value = 12.7180411142318
$ws.Range("F7").Value = 28.17752165524932
$ws.Range("G7").Value = 35.68818595479057
$ws.Range("H7").Value = 15.89958271300986
$ws.Range("I7").Value = 24.67392560865132
$ws.Range("L7").Value = 9.54946786682463
$ws.Range("M7").Value = 53.92101092238529

# Row 8
$ws.Range("C8").Value = 8.321848774512596
$ws.Range("D8").Value = 5.25345012256606
$ws.Range("E8").Value = 13.13288952245017
$ws.Range("F8").Value = 27.73469741540586
$ws.Range("G8").Value = 35.00268220307039
$ws.Range("H8").Value = 15.61508709806954
$ws.Range("I8").Value = 24.06459050756882
$ws.Range("L8").Value = 9.74412691587964
$ws.Range("M8").Value = 58.35630826692699

# Row 9
$ws.Range("C9").Value = 8.4528879681798
$ws.Range("D9").Value = 5.009822844477774
$ws.Range("E9").Value = 13.94847913057338
$ws.Range("F9").Value = 27.06989733724298
$ws.Range("G9").Value = 34.02720721065016
$ws.Range("H9").Value = 15.12677457779945
$ws.Range("I9").Value = 23.00992177172587
$ws.Range("L9").Value = 10.14429895073718
$ws.Range("M9").Value = 66.1485758021221

# Row 10
$ws.Range("C10").Value = 8.559946473421203
$ws.Range("D10").Value = 4.842165638203094
$ws.Range("E10").Value = 14.54297155420957
$ws.Range("F10").Value = 26.71775748852273
$ws.Range("G10").Value = 33.563691794396
$ws.Range("H10").Value = 14.81283284822952
$ws.Range("I10").Value = 22.32685632121979
$ws.Range("L10").Value = 10.44727233331916
$ws.Range("M10").Value = 71.31408504220791

# Row 11
$ws.Range("C11").Value = 8.610916070264178
$ws.Range("D11").Value = 4.76834317184181
$ws.Range("E11").Value = 14.81127701971433
$ws.Range("F11").Value = 26.58984282040767
$ws.Range("G11").Value = 33.41393175969338
$ws.Range("H11").Value = 14.68032469790476
$ws.Range("I11").Value = 22.03763131069893
$ws.Range("L11").Value = 10.58658467186598
$ws.Range("M11").Value = 73.54199597523962

# Row 12
$ws.Range("C12").Value = 8.630536209464397
$ws.Range("D12").Value = 4.740740348746761
$ws.Range("E12").Value = 14.91248806934514
$ws.Range("F12").Value = 26.54627070521997
$ws.Range("G12").Value = 33.36650551386546
$ws.Range("H12").Value = 14.6316822671976
$ws.Range("I12").Value = 21.93135145268167
$ws.Range("L12").Value = 10.63951559126664
$ws.Range("M12").Value = 74.36812039208836

# Row 13
$ws.Range("C13").Value = 8.62629658517281
$ws.Range("D13").Value = 4.746669433121654
$ws.Range("E13").Value = 14.89070906726062
$ws.Range("F13").Value = 26.55543471196807
$ws.Range("G13").Value = 33.37629884877116
$ws.Range("H13").Value = 14.64208909652708
$ws.Range("I13").Value = 21.95409390908066
$ws.Range("L13").Value = 10.6281087049924
$ws.Range("M13").Value = 74.19097946660473

# Row 14
$ws.Range("C14").Value = 8.61252387757866
$ws.Range("D14").Value = 4.76606520814518
$ws.Range("E14").Value = 14.81961198823064
$ws.Range("F14").Value = 26.58615929396134
$ws.Range("G14").Value = 33.40984107928497
$ws.Range("H14").Value = 14.67629176073431
$ws.Range("I14").Value = 22.02882168346342
$ws.Range("L14").Value = 10.59093600450882
$ws.Range("M14").Value = 73.61031312680805

# Row 15
$ws.Range("C15").Value = 8.604129034969169
$ws.Range("D15").Value = 4.777991566710711
$ws.Range("E15").Value = 14.77600978799218
$ws.Range("F15").Value = 26.60561943384128
$ws.Range("G15").Value = 33.43161041434498
$ws.Range("H15").Value = 14.69744350879609
$ws.Range("I15").Value = 22.07502171905356
$ws.Range("L15").Value = 10.56818853605738
$ws.Range("M15").Value = 73.25235470624031

# Row 16
$ws.Range("C16").Value = 8.556660527096289
$ws.Range("D16").Value = 4.84703934469294
$ws.Range("E16").Value = 14.52538690141111
$ws.Range("F16").Value = 26.7267857066903
$ws.Range("G16").Value = 33.57475127470885
$ws.Range("H16").Value = 14.82170443552283
$ws.Range("I16").Value = 22.34620316492775
$ws.Range("L16").Value = 10.43819462892768
$ws.Range("M16").Value = 71.1660314147687

# Row 17
$ws.Range("C17").Value = 8.528116050145487
$ws.Range("D17").Value = 4.890024399804674
$ws.Range("E17").Value = 14.37102950358765
$ws.Range("F17").Value = 26.80954367826582
$ws.Range("G17").Value = 33.67856985188855
$ws.Range("H17").Value = 14.90061091115118
$ws.Range("I17").Value = 22.51817499003436
$ws.Range("L17").Value = 10.35880135958185
$ws.Range("M17").Value = 69.8549011262899

# Row 18
$ws.Range("C18").Value = 8.51191194581155
$ws.Range("D18").Value = 4.914978391280389
$ws.Range("E18").Value = 14.28205101228759
$ws.Range("F18").Value = 26.86016820978256
$ws.Range("G18").Value = 33.74400470576393
$ws.Range("H18").Value = 14.94696239540763
$ws.Range("I18").Value = 22.61910288284367
$ws.Range("L18").Value = 10.31327867175766
$ws.Range("M18").Value = 69.08929752487192

# Row 19
$ws.Range("C19").Value = 8.506462470475784
$ws.Range("D19").Value = 4.923466906320691
$ws.Range("E19").Value = 14.25189341767137
$ws.Range("F19").Value = 26.87782139536422
$ws.Range("G19").Value = 33.76712706079969
$ws.Range("H19").Value = 14.962820649518
$ws.Range("I19").Value = 22.65361691768747
$ws.Range("L19").Value = 10.29789107693124
$ws.Range("M19").Value = 68.82810702219818

# Row 20
$ws.Range("C20").Value = 8.531132569220031
$ws.Range("D20").Value = 4.885424749359954
$ws.Range("E20").Value = 14.38748205170868
$ws.Range("F20").Value = 26.80041931559613
$ws.Range("G20").Value = 33.66692241916407
$ws.Range("H20").Value = 14.89211077013988
$ws.Range("I20").Value = 22.49965878137792
$ws.Range("L20").Value = 10.36723846824701
$ws.Range("M20").Value = 69.99566107514561

# Row 21
$ws.Range("C21").Value = 8.616560653449676
$ws.Range("D21").Value = 4.760358634150428
$ws.Range("E21").Value = 14.84050614543429
$ws.Range("F21").Value = 26.57700086803239
$ws.Range("G21").Value = 33.39973301118179
$ws.Range("H21").Value = 14.66620348860368
$ws.Range("I21").Value = 22.00678302337701
$ws.Range("L21").Value = 10.60185003137328
$ws.Range("M21").Value = 73.78134475998563

# Row 22
$ws.Range("C22").Value = 8.674249227320324
$ws.Range("D22").Value = 4.680674061023824
$ws.Range("E22").Value = 15.13427507826838
$ws.Range("F22").Value = 26.45945376938913
$ws.Range("G22").Value = 33.27944890888703
$ws.Range("H22").Value = 14.52753796871228
$ws.Range("I22").Value = 21.70364594954614
$ws.Range("L22").Value = 10.75619423674647
$ws.Range("M22").Value = 76.15330611730649

# Row 23
$ws.Range("C23").Value = 8.643292255040732
$ws.Range("D23").Value = 4.723014970144805
$ws.Range("E23").Value = 14.97772180568416
$ws.Range("F23").Value = 26.51951134809989
$ws.Range("G23").Value = 33.3385132014551
$ws.Range("H23").Value = 14.60070614580684
$ws.Range("I23").Value = 21.86364492679881
$ws.Range("L23").Value = 10.67373713581614
$ws.Range("M23").Value = 74.89669163653524

# Row 24
$ws.Range("C24").Value = 8.52976815722068
$ws.Range("D24").Value = 4.887503499721755
$ws.Range("E24").Value = 14.38004458447881
$ws.Range("F24").Value = 26.80453496980931
$ws.Range("G24").Value = 33.67217035723993
$ws.Range("H24").Value = 14.89595061330353
$ws.Range("I24").Value = 22.50802355722864
$ws.Range("L24").Value = 10.36342366990584
$ws.Range("M24").Value = 69.93206030973028

# Row 25
$ws.Range("C25").Value = 8.415519494163556
$ws.Range("D25").Value = 5.073733961296958
$ws.Range("E25").Value = 13.72826339107283
$ws.Range("F25").Value = 27.22670774926138
$ws.Range("G25").Value = 34.24871423278877
$ws.Range("H25").Value = 15.2511898234007
$ws.Range("I25").Value = 23.27961238270911
$ws.Range("L25").Value = 10.03430466887824
$ws.Range("M25").Value = 64.13848133436339

